# Marchează prezențele din săptămâna 5 (coloana G) pentru studenți.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rândurile studenților prezenți în săptămâna 5 (coloana G), la fel ca la
# celelalte săptămâni deja completate. Rândurile 5, 9 și 15 rămân goale
# (studenți absenți în săptămâna 5).
$rows = @(3, 4, 6, 7, 8, 10, 11, 12, 13, 14, 16, 17, 18, 19, 20, 21)

foreach ($r in $rows) {
    $ws.Range("G$r").Value = 1
}

# Actualizează celula activă selectată în panoul din dreapta-jos.
$ws.Range("I12").Select()

# Pragul pentru formatarea condiționată a coloanei Prezențe (Q) scade de la
# 8 la 4, reflectând numărul minim de prezențe așteptat.
$cf = $ws.Range("Q3:Q52").FormatConditions.Item(1)
$cf.Formula1 = "4"
